# Updates odds values in the "Jogos da Semana" FlashScore sheet.
# Only numeric odds/column cells in rows 2, 4, 6, 7 and 10 change; all
# other cells (Id/Date/Time/League/Home/Away, headers, other rows) stay
# the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 32 updated odds cells
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 1.83
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 5.5
$ws.Range("AA2").Value = 23
$ws.Range("AC2").Value = 6
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 81
$ws.Range("AH2").Value = 7.5
$ws.Range("AK2").Value = 41
$ws.Range("AM2").Value = 51
$ws.Range("AO2").Value = 15
$ws.Range("AS2").Value = 351
$ws.Range("AT2").Value = 2.2
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 81
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 351

# Row 4: 13 updated odds cells
$ws.Range("G4").Value = 2.05
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 3.3
$ws.Range("J4").Value = 2.63
$ws.Range("U4").Value = 1.57
$ws.Range("V4").Value = 2.25
$ws.Range("AB4").Value = 23
$ws.Range("AI4").Value = 19
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 23
$ws.Range("AS4").Value = 126
$ws.Range("AX4").Value = 17
$ws.Range("BD4").Value = 151

# Row 6: 16 updated odds cells
$ws.Range("G6").Value = 1.9
$ws.Range("I6").Value = 4.2
$ws.Range("J6").Value = 2.6
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("X6").Value = 8.5
$ws.Range("Z6").Value = 15
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 15
$ws.Range("AN6").Value = 3.75
$ws.Range("AO6").Value = 10
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 34
$ws.Range("AW6").Value = 6
$ws.Range("AX6").Value = 23
$ws.Range("AY6").Value = 34

# Row 7: 4 updated odds cells
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("O7").Value = 1.62
$ws.Range("P7").Value = 2.2

# Row 10: 18 updated odds cells
$ws.Range("G10").Value = 3
$ws.Range("I10").Value = 2.38
$ws.Range("J10").Value = 3.75
$ws.Range("L10").Value = 3.1
$ws.Range("N10").Value = 8.5
$ws.Range("W10").Value = 8.5
$ws.Range("X10").Value = 15
$ws.Range("AA10").Value = 26
$ws.Range("AD10").Value = 6
$ws.Range("AI10").Value = 11
$ws.Range("AJ10").Value = 9.5
$ws.Range("AK10").Value = 23
$ws.Range("AN10").Value = 5
$ws.Range("AP10").Value = 29
$ws.Range("AW10").Value = 4.33
$ws.Range("AX10").Value = 13
$ws.Range("AZ10").Value = 41
$ws.Range("BA10").Value = 67
